$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new time-tracking rows (90-92).
# Cell values are set in the order that reproduces the expected shared-string
# table layout (new unique strings appended as: client/api, /api/users/:id
# route..., Analyzer save nappi..., favoriteService...).
$ws.Range("B90").Value = 3
$ws.Range("D90").Value = "client/api"

$ws.Range("B91").Value = 1
$ws.Range("C91").Value = "/api/users/:id route tehty ja testattu"
$ws.Range("D91").Value = "api"

$ws.Range("C90").Value = "Analyzer save nappi ja profiilin pohjan rakenne, backend korjaus ('validoi' vahingossa käyttäjän)"

$ws.Range("B92").Value = 1
$ws.Range("C92").Value = "favoriteService, ongelmat 'populate' kohdassa API:n kanssa korjattu, tuo oikean datan, testattu"
$ws.Range("D92").Value = "client/api"

# Apply centered style to the new B cells, matching the rest of the column
$ws.Range("B90:B92").HorizontalAlignment = -4108
$ws.Range("B90:B92").VerticalAlignment = -4108

# Update the totals formula to include the new rows, then recalc
$ws.Range("B96").Formula = "=SUM(B2:B92)"

$excel.Calculate()

# Update the view: selected cell (matches the target selection)
$null = $ws.Range("C92").Select()
